# Apply the updated crypto price/volume snapshot scraped on
# Thu May 18 23:01:46 UTC 2023. Column D holds plain-text numeric-looking
# strings (e.g. "1.000", "26.909.83") that must round-trip as literal text
# rather than being auto-coerced into numbers by Excel, so each Price cell
# is stamped with a text NumberFormat before the value is written, then
# restored to the workbook default "Normal" style so no stray formatting
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ column letter = new value }
$updates = [ordered]@{
    2 = @{ 'D'='26.907.96' }
    3 = @{ 'D'='1.810.19'; 'E'='  -0.61%  ' }
    4 = @{ 'D'='1.001'; 'E'='  -0.02%  ' }
    5 = @{ 'D'='310.70'; 'E'='  -1.02%  ' }
    7 = @{ 'D'='0.4624'; 'E'='  +3.43%  ' }
    8 = @{ 'D'='0.3749'; 'E'='  -0.08%  ' }
    9 = @{ 'D'='0.07418'; 'E'='  -1.12%  ' }
    10 = @{ 'D'='0.8774'; 'E'='  -0.87%  ' }
    11 = @{ 'E'='  -2.56%  ' }
    12 = @{ 'D'='1.815.95' }
    13 = @{ 'D'='5.356'; 'E'='  -0.86%  ' }
    14 = @{ 'D'='6.548'; 'E'='  -2.99%  ' }
    15 = @{ 'D'='0.07053'; 'E'='  -0.71%  ' }
    16 = @{ 'D'='91.46'; 'E'='  -2.52%  ' }
    17 = @{ 'D'='1.002'; 'E'='  +0.11%  ' }
    18 = @{ 'D'='0.000008768'; 'E'='  -0.43%  ' }
    19 = @{ 'D'='1.000'; 'E'='  -0.01%  ' }
    20 = @{ 'D'='14.77'; 'E'='  -2.72%  ' }
    21 = @{ 'D'='26.911.25'; 'E'='  -1.64%  ' }
    22 = @{ 'D'='5.311'; 'E'='  +1.03%  ' }
    23 = @{ 'D'='10.91'; 'E'='  -0.15%  ' }
    24 = @{ 'D'='2.001.65'; 'E'='  -2.67%  ' }
    25 = @{ 'D'='1.922'; 'E'='  -2.23%  ' }
    26 = @{ 'D'='151.50'; 'E'='  -0.11%  ' }
    27 = @{ 'D'='18.61'; 'E'='  +0.17%  ' }
    28 = @{ 'D'='2.149'; 'E'='  -9.50%  ' }
    29 = @{ 'E'='  -1.05%  ' }
    30 = @{ 'D'='117.32'; 'E'='  -0.45%  ' }
    31 = @{ 'D'='0.08892'; 'E'='  +0.69%  ' }
    32 = @{ 'D'='0.7714'; 'E'='  -2.62%  ' }
    33 = @{ 'D'='1.168'; 'E'='  -2.48%  ' }
    34 = @{ 'D'='4.485'; 'E'='  -0.42%  ' }
    35 = @{ 'D'='2.897'; 'E'='  -0.45%  ' }
    36 = @{ 'D'='1.000'; 'E'='  +0.00%  ' }
    37 = @{ 'D'='1.120'; 'E'='  +0.57%  ' }
    38 = @{ 'D'='0.01957'; 'E'='  -1.91%  ' }
    39 = @{ 'D'='2.435'; 'E'='  +4.25%  ' }
    40 = @{ 'D'='0.05244'; 'E'='  -1.60%  ' }
    41 = @{ 'B'='TheSandbox'; 'C'='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; 'D'='0.5349'; 'E'='  +0.64%  ' }
    42 = @{ 'B'='FraxShare'; 'C'='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; 'D'='7.263'; 'E'='  -1.51%  ' }
    43 = @{ 'D'='2.906'; 'E'='  +1.71%  ' }
    44 = @{ 'D'='0.1662'; 'E'='  -3.54%  ' }
    45 = @{ 'D'='8.562'; 'E'='  -1.96%  ' }
    46 = @{ 'D'='0.5069'; 'E'='  -0.73%  ' }
    47 = @{ 'D'='10.34'; 'E'='  -3.46%  ' }
    48 = @{ 'D'='104.20'; 'E'='  -1.49%  ' }
    49 = @{ 'D'='1.0000'; 'E'='  -0.02%  ' }
    50 = @{ 'D'='1.666'; 'E'='  -2.00%  ' }
    51 = @{ 'D'='0.06329'; 'E'='  -0.78%  ' }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.Contains('D')) {
        $ws.Range("D$row").NumberFormat = "@"
    }
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
    if ($rowData.Contains('D')) {
        $ws.Range("D$row").Style = 'Normal'
    }
}
